$d = $word.ActiveDocument

$replacements = @(
    @("388÷7=55, 3", "881÷2=440, 1"),
    @("265÷8=33, 1", "148÷5=29, 3"),
    @("383÷9=42, 5", "436÷5=87, 1"),
    @("162÷8=20, 2", "640÷4=160, 0"),
    @("733÷5=146, 3", "918÷4=229, 2"),
    @("153÷6=25, 3", "652÷4=163, 0"),
    @("322÷9=35, 7", "915÷2=457, 1"),
    @("233÷4=58, 1", "494÷3=164, 2"),
    @("809÷5=161, 4", "320÷3=106, 2"),
    @("400÷9=44, 4", "212÷7=30, 2"),
    @("558÷5=111, 3", "901÷5=180, 1"),
    @("936÷9=104, 0", "672÷8=84, 0"),
    @("198÷2=99, 0", "361÷9=40, 1"),
    @("228÷7=32, 4", "538÷2=269, 0"),
    @("483÷7=69, 0", "700÷6=116, 4"),
    @("525÷3=175, 0", "302÷4=75, 2"),
    @("306÷9=34, 0", "289÷5=57, 4"),
    @("216÷9=24, 0", "487÷7=69, 4"),
    @("865÷9=96, 1", "201÷8=25, 1"),
    @("464÷5=92, 4", "978÷9=108, 6"),
    @("427÷6=71, 1", "379÷9=42, 1"),
    @("254÷5=50, 4", "388÷5=77, 3"),
    @("977÷6=162, 5", "733÷9=81, 4"),
    @("485÷3=161, 2", "771÷7=110, 1"),
    @("479÷4=119, 3", "717÷7=102, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
